$d = $word.ActiveDocument

# 1. Fix the double space in the "Ejemplo..." paragraph:
#    "Ejemplo de un documento de  pase o reporte de código "
# -> "Ejemplo de un documento de pase o reporte de código "
$d.Content.Find.Execute(
    "de  pase", $true, $false, $false, $false, $false,
    $true, 1, $false, "de pase", 2
)

# 2. Append a new blank paragraph after the last paragraph.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$blankXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@
$d.Paragraphs.Last.Range.InsertXML($blankXml)

# 3. Append the "Segundo commit" paragraph (with proofErr spell-check markers
#    around the English word "commit").
$d.Paragraphs.Last.Range.InsertParagraphAfter()

$segundoXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">Segundo </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>commit</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
"@
$d.Paragraphs.Last.Range.InsertXML($segundoXml)

# 4. Append the final "Agregando una observación..." paragraph.
$d.Paragraphs.Last.Range.InsertParagraphAfter()

$obsXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>Agregando una observación  nueva no comentada anteriormente a  última hora</w:t>
  </w:r>
</w:p>
"@
$d.Paragraphs.Last.Range.InsertXML($obsXml)

Write-Output "edit complete"
